$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# 477×2= -> 599×3=
$cell = $t.Cell(1, 1)
$cell.Range.Text = "599×3="

# 310×7= -> 613×4=
$cell = $t.Cell(1, 2)
$cell.Range.Text = "613×4="

# 360×3= -> 307×3=
$cell = $t.Cell(1, 3)
$cell.Range.Text = "307×3="

# 845×4= -> 771×8=
$cell = $t.Cell(1, 4)
$cell.Range.Text = "771×8="

# 700×5= -> 762×4=
$cell = $t.Cell(1, 5)
$cell.Range.Text = "762×4="

# 895×4= -> 317×9=
$cell = $t.Cell(5, 1)
$cell.Range.Text = "317×9="

# 242×4= -> 119×7=
$cell = $t.Cell(5, 2)
$cell.Range.Text = "119×7="

# 125×4= -> 537×7=
$cell = $t.Cell(5, 3)
$cell.Range.Text = "537×7="

# 182×7= -> 887×2=
$cell = $t.Cell(5, 4)
$cell.Range.Text = "887×2="

# 398×2= -> 453×3=
$cell = $t.Cell(5, 5)
$cell.Range.Text = "453×3="

# 675×6= -> 638×7=
$cell = $t.Cell(10, 1)
$cell.Range.Text = "638×7="

# 346×9= -> 749×5=
$cell = $t.Cell(10, 2)
$cell.Range.Text = "749×5="

# 332×2= -> 263×7=
$cell = $t.Cell(10, 3)
$cell.Range.Text = "263×7="

# 632×8= -> 255×2=
$cell = $t.Cell(10, 4)
$cell.Range.Text = "255×2="

# 395×9= -> 632×9=
$cell = $t.Cell(10, 5)
$cell.Range.Text = "632×9="

# 637×3= -> 295×3=
$cell = $t.Cell(15, 1)
$cell.Range.Text = "295×3="

# 310×7= -> 839×8=
$cell = $t.Cell(15, 2)
$cell.Range.Text = "839×8="

# 817×8= -> 859×5=
$cell = $t.Cell(15, 3)
$cell.Range.Text = "859×5="

# 239×7= -> 483×9=
$cell = $t.Cell(15, 4)
$cell.Range.Text = "483×9="

# 920×4= -> 723×9=
$cell = $t.Cell(15, 5)
$cell.Range.Text = "723×9="

# 218×6= -> 850×3=
$cell = $t.Cell(20, 1)
$cell.Range.Text = "850×3="

# 326×4= -> 803×3=
$cell = $t.Cell(20, 2)
$cell.Range.Text = "803×3="

# 214×8= -> 876×3=
$cell = $t.Cell(20, 3)
$cell.Range.Text = "876×3="

# 932×4= -> 954×4=
$cell = $t.Cell(20, 4)
$cell.Range.Text = "954×4="

# 494×3= -> 578×2=
$cell = $t.Cell(20, 5)
$cell.Range.Text = "578×2="

